$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": header refresh + updated schedule rows + one new row
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:28:22"
$ws1.Range("A3").Value = "Total filas: 3"

# Row 6: 14_ABASTO -> 215_ALUAR
$ws1.Range("A6").Value = "02:28:22"
$ws1.Range("B6").Value = "02:58"
$ws1.Range("C6").Value = "215_ALUAR"
$ws1.Range("D6").Value = 30
$ws1.Range("E6").Value = "LP1912"

# Row 7: 215_ALUAR -> 14_ABASTO
$ws1.Range("A7").Value = "02:28:22"
$ws1.Range("B7").Value = "03:48"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 80
$ws1.Range("E7").Value = "LP1912"

# Row 8: new row - 81_EL PELIGRO
$ws1.Range("A8").Value = "02:28:22"
$ws1.Range("B8").Value = "04:01"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 93
$ws1.Range("E8").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": header refresh + updated minutes on existing row
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:28:22"
$ws2.Range("A6").Value = "02:28:22"
$ws2.Range("D6").Value = 30

# ---------------------------------------------------------------------------
# Sheet "6203-6173": header refresh only
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:28:22"
